$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style source cells (existing rows) used to replicate formatting without creating new style entries
$dateFmtSrc = $ws.Range("B454")      # numFmtId 14 (date) -> style index 1
$centerFmtSrc = $ws.Range("D454")    # centered alignment -> style index 3
$namePosFmtSrc = $ws.Range("E245")   # style index 6 (used for E/F columns in a couple of rows)

function Set-StyledCell($cellRef, $styleSrc, $value) {
    $target = $ws.Range($cellRef)
    $styleSrc.Copy()
    $target.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $target.Value = $value
}

# ---- Row 455 ----
$ws.Range("A455").Value = "Entrainement"
Set-StyledCell "B455" $dateFmtSrc 45908
$ws.Range("C455").Value = "Global"
$ws.Range("E455").Value = "Karahali Souaré"
$ws.Range("F455").Value = "right forward"
$ws.Range("G455").Value = "01:04:47"
$ws.Range("H455").Value = 7.73
$ws.Range("I455").Value = 0.57999999999999996
$ws.Range("J455").Value = 7.13
$ws.Range("K455").Value = 0.56000000000000005
$ws.Range("L455").Value = 0.04
$ws.Range("M455").Value = 0
$ws.Range("N455").Value = 0
$ws.Range("O455").Value = 1
$ws.Range("P455").Value = 6.85
$ws.Range("Q455").Value = 25.02
$ws.Range("R455").Value = 5.27
$ws.Range("S455").Value = 40
$ws.Range("T455").Value = 18
$ws.Range("U455").Value = 19
$ws.Range("V455").Value = 16

# ---- Row 456 ----
$ws.Range("A456").Value = "Entrainement"
Set-StyledCell "B456" $dateFmtSrc 45908
$ws.Range("C456").Value = "Global"
$ws.Range("E456").Value = "Omar Benyounes"
$ws.Range("F456").Value = "center midfield"
$ws.Range("G456").Value = "01:02:45"
$ws.Range("H456").Value = 7.74
$ws.Range("I456").Value = 0.56000000000000005
$ws.Range("J456").Value = 7.17
$ws.Range("K456").Value = 0.53
$ws.Range("L456").Value = 0.03
$ws.Range("M456").Value = 0
$ws.Range("N456").Value = 0
$ws.Range("O456").Value = 1
$ws.Range("P456").Value = 7.36
$ws.Range("Q456").Value = 25.29
$ws.Range("R456").Value = 4.5199999999999996
$ws.Range("S456").Value = 26
$ws.Range("T456").Value = 4
$ws.Range("U456").Value = 13
$ws.Range("V456").Value = 5

# ---- Row 457 ----
$ws.Range("A457").Value = "Entrainement"
Set-StyledCell "B457" $dateFmtSrc 45908
$ws.Range("C457").Value = "Global"
$ws.Range("E457").Value = "Ilyes Boughanmi"
$ws.Range("F457").Value = "center forward"
$ws.Range("G457").Value = "01:03:35"
$ws.Range("H457").Value = 7.18
$ws.Range("I457").Value = 0.17
$ws.Range("J457").Value = 7.01
$ws.Range("K457").Value = 0.17
$ws.Range("L457").Value = 0
$ws.Range("M457").Value = 0
$ws.Range("N457").Value = 0
$ws.Range("O457").Value = 0
$ws.Range("P457").Value = 6.74
$ws.Range("Q457").Value = 19.29
$ws.Range("R457").Value = 4.49
$ws.Range("S457").Value = 24
$ws.Range("T457").Value = 4
$ws.Range("U457").Value = 12
$ws.Range("V457").Value = 1

# ---- Row 458 ----
$ws.Range("A458").Value = "Entrainement"
Set-StyledCell "B458" $dateFmtSrc 45908
$ws.Range("C458").Value = "Global"
$ws.Range("E458").Value = "Mattheo Haon"
$ws.Range("F458").Value = "right back"
$ws.Range("G458").Value = "01:04:18"
$ws.Range("H458").Value = 7.42
$ws.Range("I458").Value = 0.49
$ws.Range("J458").Value = 6.93
$ws.Range("K458").Value = 0.49
$ws.Range("L458").Value = 0
$ws.Range("M458").Value = 0
$ws.Range("N458").Value = 0
$ws.Range("O458").Value = 0
$ws.Range("P458").Value = 6.89
$ws.Range("Q458").Value = 19.100000000000001
$ws.Range("R458").Value = 4.2
$ws.Range("S458").Value = 26
$ws.Range("T458").Value = 1
$ws.Range("U458").Value = 14
$ws.Range("V458").Value = 2

# ---- Row 459 ----
$ws.Range("A459").Value = "Entrainement"
Set-StyledCell "B459" $dateFmtSrc 45908
$ws.Range("C459").Value = "Global"
$ws.Range("E459").Value = "Ilan Ihaddadene"
$ws.Range("F459").Value = "center midfield"
$ws.Range("G459").Value = "01:04:18"
$ws.Range("H459").Value = 8.27
$ws.Range("I459").Value = 0.39
$ws.Range("J459").Value = 7.88
$ws.Range("K459").Value = 0.39
$ws.Range("L459").Value = 0
$ws.Range("M459").Value = 0
$ws.Range("N459").Value = 0
$ws.Range("O459").Value = 0
$ws.Range("P459").Value = 7.66
$ws.Range("Q459").Value = 18.71
$ws.Range("R459").Value = 4.0999999999999996
$ws.Range("S459").Value = 20
$ws.Range("T459").Value = 1
$ws.Range("U459").Value = 7
$ws.Range("V459").Value = 1

# ---- Row 460 ----
$ws.Range("A460").Value = "Entrainement"
Set-StyledCell "B460" $dateFmtSrc 45908
$ws.Range("C460").Value = "Global"
$ws.Range("E460").Value = "Levy Ndoutoume"
$ws.Range("F460").Value = "left back"
$ws.Range("G460").Value = "01:02:45"
$ws.Range("H460").Value = 7.11
$ws.Range("I460").Value = 0.11
$ws.Range("J460").Value = 6.99
$ws.Range("K460").Value = 0.11
$ws.Range("L460").Value = 0.01
$ws.Range("M460").Value = 0
$ws.Range("N460").Value = 0
$ws.Range("O460").Value = 0
$ws.Range("P460").Value = 6.02
$ws.Range("Q460").Value = 21.46
$ws.Range("R460").Value = 4.41
$ws.Range("S460").Value = 21
$ws.Range("T460").Value = 4
$ws.Range("U460").Value = 17
$ws.Range("V460").Value = 3

# ---- Row 461 ----
$ws.Range("A461").Value = "Entrainement"
Set-StyledCell "B461" $dateFmtSrc 45909
$ws.Range("C461").Value = "Global"
Set-StyledCell "D461" $centerFmtSrc "J+3"
$ws.Range("E461").Value = "Yanis Berrached"
$ws.Range("F461").Value = "center midfield"
$ws.Range("G461").Value = "01:22:10"
$ws.Range("H461").Value = 6.04
$ws.Range("I461").Value = 0.15
$ws.Range("J461").Value = 5.88
$ws.Range("K461").Value = 0.13
$ws.Range("L461").Value = 0.02
$ws.Range("M461").Value = 0
$ws.Range("N461").Value = 0
$ws.Range("O461").Value = 0
$ws.Range("P461").Value = 4.45
$ws.Range("Q461").Value = 23.35
$ws.Range("R461").Value = 3.98
$ws.Range("S461").Value = 19
$ws.Range("T461").Value = 0
$ws.Range("U461").Value = 12
$ws.Range("V461").Value = 2

# ---- Row 462 ----
$ws.Range("A462").Value = "Entrainement"
Set-StyledCell "B462" $dateFmtSrc 45909
$ws.Range("C462").Value = "Global"
Set-StyledCell "D462" $centerFmtSrc "J+3"
$ws.Range("E462").Value = "Mattheo Haon"
$ws.Range("F462").Value = "right back"
$ws.Range("G462").Value = "01:31:13"
$ws.Range("H462").Value = 6.98
$ws.Range("I462").Value = 0.32
$ws.Range("J462").Value = 6.66
$ws.Range("K462").Value = 0.31
$ws.Range("L462").Value = 0.02
$ws.Range("M462").Value = 0
$ws.Range("N462").Value = 0
$ws.Range("O462").Value = 0
$ws.Range("P462").Value = 4.53
$ws.Range("Q462").Value = 22.28
$ws.Range("R462").Value = 4.1399999999999997
$ws.Range("S462").Value = 48
$ws.Range("T462").Value = 3
$ws.Range("U462").Value = 23
$ws.Range("V462").Value = 4

# ---- Row 463 ----
$ws.Range("A463").Value = "Entrainement"
Set-StyledCell "B463" $dateFmtSrc 45909
$ws.Range("C463").Value = "Global"
Set-StyledCell "D463" $centerFmtSrc "J+3"
$ws.Range("E463").Value = "Ilyes Boughanmi"
$ws.Range("F463").Value = "center forward"
$ws.Range("G463").Value = "01:31:21"
$ws.Range("H463").Value = 5.63
$ws.Range("I463").Value = 0.16
$ws.Range("J463").Value = 5.46
$ws.Range("K463").Value = 0.17
$ws.Range("L463").Value = 0
$ws.Range("M463").Value = 0
$ws.Range("N463").Value = 0
$ws.Range("O463").Value = 0
$ws.Range("P463").Value = 3.61
$ws.Range("Q463").Value = 20.100000000000001
$ws.Range("R463").Value = 4.05
$ws.Range("S463").Value = 27
$ws.Range("T463").Value = 1
$ws.Range("U463").Value = 21
$ws.Range("V463").Value = 2

# ---- Row 464 ----
$ws.Range("A464").Value = "Entrainement"
Set-StyledCell "B464" $dateFmtSrc 45909
$ws.Range("C464").Value = "Global"
Set-StyledCell "D464" $centerFmtSrc "J+3"
$ws.Range("E464").Value = "Ilan Ihaddadene"
$ws.Range("F464").Value = "center midfield"
$ws.Range("G464").Value = "01:31:48"
$ws.Range("H464").Value = 7.24
$ws.Range("I464").Value = 0.3
$ws.Range("J464").Value = 6.94
$ws.Range("K464").Value = 0.28000000000000003
$ws.Range("L464").Value = 0.03
$ws.Range("M464").Value = 0
$ws.Range("N464").Value = 0
$ws.Range("O464").Value = 0
$ws.Range("P464").Value = 4.6500000000000004
$ws.Range("Q464").Value = 21.82
$ws.Range("R464").Value = 4.33
$ws.Range("S464").Value = 50
$ws.Range("T464").Value = 3
$ws.Range("U464").Value = 22
$ws.Range("V464").Value = 4

# ---- Row 465 ----
$ws.Range("A465").Value = "Entrainement"
Set-StyledCell "B465" $dateFmtSrc 45909
$ws.Range("C465").Value = "Global"
Set-StyledCell "D465" $centerFmtSrc "J+3"
$ws.Range("E465").Value = "Jeremie Laurent"
$ws.Range("F465").Value = "left forward"
$ws.Range("G465").Value = "01:26:15"
$ws.Range("H465").Value = 6.68
$ws.Range("I465").Value = 0.32
$ws.Range("J465").Value = 6.35
$ws.Range("K465").Value = 0.31
$ws.Range("L465").Value = 0.02
$ws.Range("M465").Value = 0
$ws.Range("N465").Value = 0
$ws.Range("O465").Value = 0
$ws.Range("P465").Value = 4.59
$ws.Range("Q465").Value = 22.54
$ws.Range("R465").Value = 4.26
$ws.Range("S465").Value = 42
$ws.Range("T465").Value = 4
$ws.Range("U465").Value = 36
$ws.Range("V465").Value = 4

# ---- Row 466 ----
$ws.Range("A466").Value = "Entrainement"
Set-StyledCell "B466" $dateFmtSrc 45909
$ws.Range("C466").Value = "Global"
Set-StyledCell "D466" $centerFmtSrc "J+3"
$ws.Range("E466").Value = "Omar Benyounes"
$ws.Range("F466").Value = "center midfield"
$ws.Range("G466").Value = "01:31:30"
$ws.Range("H466").Value = 6.69
$ws.Range("I466").Value = 0.36
$ws.Range("J466").Value = 6.32
$ws.Range("K466").Value = 0.34
$ws.Range("L466").Value = 0.03
$ws.Range("M466").Value = 0
$ws.Range("N466").Value = 0
$ws.Range("O466").Value = 0
$ws.Range("P466").Value = 4.3
$ws.Range("Q466").Value = 24.03
$ws.Range("R466").Value = 4.6900000000000004
$ws.Range("S466").Value = 43
$ws.Range("T466").Value = 7
$ws.Range("U466").Value = 33
$ws.Range("V466").Value = 6

# ---- Row 467 ----
$ws.Range("A467").Value = "Entrainement"
Set-StyledCell "B467" $dateFmtSrc 45909
$ws.Range("C467").Value = "Global"
Set-StyledCell "D467" $centerFmtSrc "J+3"
$ws.Range("E467").Value = "Hedi Nasri"
$ws.Range("F467").Value = "right back"
$ws.Range("G467").Value = "01:31:40"
$ws.Range("H467").Value = 5.96
$ws.Range("I467").Value = 0.19
$ws.Range("J467").Value = 5.77
$ws.Range("K467").Value = 0.18
$ws.Range("L467").Value = 0.02
$ws.Range("M467").Value = 0
$ws.Range("N467").Value = 0
$ws.Range("O467").Value = 0
$ws.Range("P467").Value = 3.88
$ws.Range("Q467").Value = 22.22
$ws.Range("R467").Value = 4.2
$ws.Range("S467").Value = 30
$ws.Range("T467").Value = 2
$ws.Range("U467").Value = 21
$ws.Range("V467").Value = 6

# ---- Row 468 ----
$ws.Range("A468").Value = "Entrainement"
Set-StyledCell "B468" $dateFmtSrc 45909
$ws.Range("C468").Value = "Global"
Set-StyledCell "D468" $centerFmtSrc "J+3"
Set-StyledCell "E468" $namePosFmtSrc "Malik Boussaid"
Set-StyledCell "F468" $namePosFmtSrc "right back"
$ws.Range("G468").Value = "01:25:40"
$ws.Range("H468").Value = 6.65
$ws.Range("I468").Value = 0.21
$ws.Range("J468").Value = 6.44
$ws.Range("K468").Value = 0.21
$ws.Range("L468").Value = 0.01
$ws.Range("M468").Value = 0
$ws.Range("N468").Value = 0
$ws.Range("O468").Value = 0
$ws.Range("P468").Value = 4.54
$ws.Range("Q468").Value = 22.34
$ws.Range("R468").Value = 4.4000000000000004
$ws.Range("S468").Value = 17
$ws.Range("T468").Value = 2
$ws.Range("U468").Value = 22
$ws.Range("V468").Value = 1

# Update selection to mirror the authored view state (top-left scroll position is not preserved by this runtime)
$excel.Goto($ws.Range("D475"), $true)

Write-Host "Added rows 455-468"
